# Update the "answers-of-two-digit_number_divided_by_one-digit_number"
# document: change the date and all 25 division answer cells.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-06 Monday" "2024-05-07 Tuesday"

Replace-Text "15÷7=2, 1" "98÷5=19, 3"
Replace-Text "23÷7=3, 2" "62÷3=20, 2"
Replace-Text "21÷5=4, 1" "29÷2=14, 1"
Replace-Text "18÷3=6, 0" "16÷5=3, 1"
Replace-Text "28÷8=3, 4" "84÷3=28, 0"

Replace-Text "26÷2=13, 0" "63÷3=21, 0"
Replace-Text "89÷6=14, 5" "64÷9=7, 1"
Replace-Text "88÷8=11, 0" "32÷3=10, 2"
Replace-Text "20÷5=4, 0" "32÷6=5, 2"
Replace-Text "24÷7=3, 3" "29÷4=7, 1"

Replace-Text "56÷9=6, 2" "36÷6=6, 0"
Replace-Text "59÷4=14, 3" "28÷4=7, 0"
Replace-Text "80÷8=10, 0" "88÷3=29, 1"
Replace-Text "44÷5=8, 4" "80÷6=13, 2"
Replace-Text "28÷6=4, 4" "67÷9=7, 4"

Replace-Text "26÷4=6, 2" "35÷8=4, 3"
Replace-Text "98÷6=16, 2" "59÷7=8, 3"
Replace-Text "62÷6=10, 2" "82÷4=20, 2"
Replace-Text "18÷4=4, 2" "99÷5=19, 4"
Replace-Text "14÷3=4, 2" "28÷9=3, 1"

Replace-Text "35÷7=5, 0" "58÷4=14, 2"
Replace-Text "82÷6=13, 4" "71÷6=11, 5"
Replace-Text "56÷8=7, 0" "39÷8=4, 7"
Replace-Text "62÷5=12, 2" "23÷9=2, 5"
Replace-Text "57÷2=28, 1" "31÷5=6, 1"
